$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows before the current row 9, pushing the footer  ---
# --- block (old rows 9-10, "total power consumption..."/"Power budget" ---
# --- summary) down to rows 11-12, leaving row 10 blank just like the   ---
# --- original layout had a blank row 8 before the footer.             ---
$ws.Range("A9:A10").EntireRow.Insert()

# --- New row 8: batteries ---
$ws.Range("B8").Value = "batteries"
$ws.Range("A8").Value = "turnigy nano-tech 180mAh 2s"
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = "NA"
$ws.Range("E8").Value = "NA"
$ws.Range("F8").Value = "NA"

# --- New header cells for mass and dimensions ---
$ws.Range("G1").Value = "mass (g)"
$ws.Range("H1").Value = "dimensions (cm lxwxh)"

$ws.Range("G8").Value = 13
$ws.Range("H8").Value = "3.5x2.0x1.0"

# --- New row 9: foam ---
$ws.Range("A9").Value = "foam"
$ws.Range("B9").Value = "foam"
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("F9").Value = "NA"
$ws.Range("G9").Value = 119
$ws.Range("H9").Value = "45x45x2.5"

# --- Move the "volume of foam" / "weight of foam" helper columns out of ---
# --- H/I (now reused by the mass/dimensions headers) into J/K.          ---
$ws.Range("J1").Value = "volume of the foam in cm^3"
$ws.Range("K1").Value = "weight of the foam in kg"
$ws.Range("J2").Formula = "=PI() * 2.5 * (45/2)^2"
$ws.Range("K2").Formula = "=(J2/1000000)*30"
$ws.Range("H2").ClearContents()
$ws.Range("I1:I2").ClearContents()

# --- Column A needs to widen to fit the new, longer part-number text ---
$ws.Columns("A").ColumnWidth = 26.75

# --- Selection moves to G2 in the saved file ---
$ws.Range("G2").Select()
